$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) "The students have 48 hours ... this document." paragraph:
#    split "document." into "do" | _GoBack bookmark | "cument.", then
#    append the new "Note: ..." sentence (with two bold spans) to the
#    same paragraph.
# ------------------------------------------------------------------
$r = $d.Content
$r.Find.Execute("this document.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$splitPos = $r.Start + 7   # right after "this do"
$bm = $d.Range($splitPos, $splitPos)
$d.Bookmarks.Add("_GoBack", $bm)

$r = $d.Content
$r.Find.Execute("this document.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$r.Collapse(0)
$r.InsertAfter(" ")
$r.Collapse(0)
$r.InsertAfter("Note: as a general rule, usually there is no deadline extension on this type of exams. And, even if administration grants an extension (e.g., for medical reasons), it should be no more than 50% of the original amount (i.e., a total of 72 hours in a 48 hour exam). If for any reason you got granted an extension longer than that, you must contact administration to verify the course responsible had agreed on such extension (there were cases in the past in which such unauthorized extensions were given by mistake). Do ")
$r.Collapse(0)
$r.InsertAfter("NOT")
$r.Font.Bold = 1
$r.Collapse(0)
$r.InsertAfter(" contact the course responsible directly, as exams must be marked anonymously. To make the exam conditions fair to all students, submissions with long extensions that were not authorized by the course responsible will be automatically evaluated as failed (i.e., an ")
$r.Font.Bold = 0
$r.Collapse(0)
$r.InsertAfter("F")
$r.Font.Bold = 1
$r.Collapse(0)
$r.InsertAfter(").")
$r.Font.Bold = 0
